$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 9 first, so the shared-string table allocates new entries
# in the same order as the target workbook (1/2, test teacher, Cô Nguyên).
$ws.Cells.Item(9,1).Value = "abc"
$ws.Cells.Item(9,2).Value = "'1/2"
$ws.Cells.Item(9,3).Value = "test teacher"
$ws.Cells.Item(9,4).Value = "Cô Nguyên"

# Row 1: update teacher date/name and add the new "Cô Nguyên" column
$ws.Cells.Item(1,2).Value = "'3/2"
$ws.Cells.Item(1,3).Value = "test teacher"
$ws.Cells.Item(1,4).Value = "Cô Nguyên"

# Row 4 is removed entirely from the sheet
$ws.Range("A4:C4").Clear()

# Row 7: update teacher date/name and add the new "Cô Nguyên" column
$ws.Cells.Item(7,2).Value = "'3/1"
$ws.Cells.Item(7,3).Value = "test teacher"
$ws.Cells.Item(7,4).Value = "Cô Nguyên"
